# Iron & Steel Foundries workbook update
# - Compact each data column (rows 2-34) upward to remove the
#   row gaps left by staggered BLS annual series imports, so
#   every series lines up on the same year row (2 = 1958 ... 18 = 1974).
# - Drop the now-empty trailing rows (19-34).
# - Restore the working selection to I6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1958
$ws.Range("B2").Value = 192.4
$ws.Range("C2").Value = "n/a"
$ws.Range("D2").Value = 162.7
$ws.Range("E2").Value = 86.86
$ws.Range("F2").Value = 2.31
$ws.Range("G2").Value = 37.6
$ws.Range("H2").Value = 1.5
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 0.6
$ws.Range("M2").Value = 2.8

# Row 3
$ws.Range("A3").Value = 1959
$ws.Range("B3").Value = 211.8
$ws.Range("C3").Value = 9.6
$ws.Range("D3").Value = 181.5
$ws.Range("E3").Value = 97.04
$ws.Range("F3").Value = 2.42
$ws.Range("G3").Value = 40.1
$ws.Range("H3").Value = 2.7
$ws.Range("I3").Value = 4.2
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 3.5
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 1.7

# Row 4
$ws.Range("A4").Value = 1960
$ws.Range("B4").Value = 204.7
$ws.Range("C4").Value = 9.5
$ws.Range("D4").Value = 173.3
$ws.Range("E4").Value = 96.61
$ws.Range("F4").Value = 2.49
$ws.Range("G4").Value = 38.8
$ws.Range("H4").Value = 2.1
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 1.2
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 0.8
$ws.Range("M4").Value = 2.7

# Row 5
$ws.Range("A5").Value = 1961
$ws.Range("B5").Value = 186.7
$ws.Range("C5").Value = 8.8
$ws.Range("D5").Value = 156.6
$ws.Range("E5").Value = 98.81
$ws.Range("F5").Value = 2.54
$ws.Range("G5").Value = 38.9
$ws.Range("H5").Value = 2.1
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 1.4
$ws.Range("K5").Value = 3.4
$ws.Range("L5").Value = 0.7
$ws.Range("M5").Value = 2

# Row 6
$ws.Range("A6").Value = 1962
$ws.Range("B6").Value = 193.6
$ws.Range("C6").Value = 8.8
$ws.Range("D6").Value = 163.7
$ws.Range("E6").Value = 106.52
$ws.Range("F6").Value = 2.63
$ws.Range("G6").Value = 40.5
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 3.5
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 3.2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.4

# Row 7
$ws.Range("A7").Value = 1963
$ws.Range("B7").Value = 198.1
$ws.Range("C7").Value = 8.8
$ws.Range("D7").Value = 168.2
$ws.Range("E7").Value = 113.28
$ws.Range("F7").Value = 2.71
$ws.Range("G7").Value = 41.8
$ws.Range("H7").Value = 3.7
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.1
$ws.Range("K7").Value = 3.2
$ws.Range("L7").Value = 1.1
$ws.Range("M7").Value = 1.2

# Row 8
$ws.Range("A8").Value = 1964
$ws.Range("B8").Value = 212.3
$ws.Range("C8").Value = 8.9
$ws.Range("D8").Value = 181.9
$ws.Range("E8").Value = 119.84
$ws.Range("F8").Value = 2.8
$ws.Range("G8").Value = 42.8
$ws.Range("H8").Value = 4.7
$ws.Range("I8").Value = 3.8
$ws.Range("J8").Value = 2.7
$ws.Range("K8").Value = 3.3
$ws.Range("L8").Value = 1.4
$ws.Range("M8").Value = 0.9

# Row 9
$ws.Range("A9").Value = 1965
$ws.Range("B9").Value = 227
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 194.6
$ws.Range("E9").Value = 125.72
$ws.Range("F9").Value = 2.89
$ws.Range("G9").Value = 43.5
$ws.Range("H9").Value = 5.5
$ws.Range("I9").Value = 4.1
$ws.Range("J9").Value = 3.2
$ws.Range("K9").Value = 3.6
$ws.Range("L9").Value = 1.8
$ws.Range("M9").Value = 0.7

# Row 10
$ws.Range("A10").Value = 1966
$ws.Range("B10").Value = 240.8
$ws.Range("C10").Value = 11.6
$ws.Range("D10").Value = 205.8
$ws.Range("E10").Value = 128.57
$ws.Range("F10").Value = 2.99
$ws.Range("G10").Value = 43
$ws.Range("H10").Value = 5.3
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 4.8
$ws.Range("L10").Value = 2.8
$ws.Range("M10").Value = 0.9

# Row 11
$ws.Range("A11").Value = 1967
$ws.Range("B11").Value = 231.3
$ws.Range("C11").Value = 12.2
$ws.Range("D11").Value = 195.3
$ws.Range("E11").Value = 127.71
$ws.Range("F11").Value = 3.07
$ws.Range("G11").Value = 41.6
$ws.Range("H11").Value = 4.2
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 4.7
$ws.Range("L11").Value = 2.4
$ws.Range("M11").Value = 1.1

# Row 12
$ws.Range("A12").Value = 1968
$ws.Range("B12").Value = 225.7
$ws.Range("C12").Value = 11.5
$ws.Range("D12").Value = 189
$ws.Range("E12").Value = 139.86
$ws.Range("F12").Value = 3.33
$ws.Range("G12").Value = 42
$ws.Range("H12").Value = 4.9
$ws.Range("I12").Value = 4.8
$ws.Range("J12").Value = 3.7
$ws.Range("K12").Value = 4.8
$ws.Range("L12").Value = 2.7
$ws.Range("M12").Value = 0.9

# Row 13
$ws.Range("A13").Value = 1969
$ws.Range("B13").Value = 235.1
$ws.Range("C13").Value = 11.7
$ws.Range("D13").Value = 197.6
$ws.Range("E13").Value = 150.23
$ws.Range("F13").Value = 3.56
$ws.Range("G13").Value = 42.2
$ws.Range("H13").Value = 5.1
$ws.Range("I13").Value = 5.8
$ws.Range("J13").Value = 4.9
$ws.Range("K13").Value = 5.6
$ws.Range("L13").Value = 3.4
$ws.Range("M13").Value = 0.7

# Row 14
$ws.Range("A14").Value = 1970
$ws.Range("B14").Value = 228.9
$ws.Range("C14").Value = 11.9
$ws.Range("D14").Value = 190.2
$ws.Range("E14").Value = 151.03
$ws.Range("F14").Value = 3.72
$ws.Range("G14").Value = 40.6
$ws.Range("H14").Value = 3.9
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 5.6
$ws.Range("L14").Value = 2.4
$ws.Range("M14").Value = 2

# Row 15
$ws.Range("A15").Value = 1971
$ws.Range("B15").Value = 217.9
$ws.Range("C15").Value = 10.7
$ws.Range("D15").Value = 180.2
$ws.Range("E15").Value = 164.43
$ws.Range("F15").Value = 4.04
$ws.Range("G15").Value = 40.7
$ws.Range("H15").Value = 3.7
$ws.Range("I15").Value = 3.6
$ws.Range("J15").Value = 2.2
$ws.Range("K15").Value = 3.9
$ws.Range("L15").Value = 1.5
$ws.Range("M15").Value = 1.4

# Row 16
$ws.Range("A16").Value = 1972
$ws.Range("B16").Value = 222.1
$ws.Range("C16").Value = 10.2
$ws.Range("D16").Value = 184.9
$ws.Range("E16").Value = 183.99
$ws.Range("F16").Value = 4.36
$ws.Range("G16").Value = 42.2
$ws.Range("H16").Value = 4.8
$ws.Range("I16").Value = 4.4
$ws.Range("J16").Value = 3.2
$ws.Range("K16").Value = 3.9
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 0.8

# Row 17
$ws.Range("A17").Value = 1973
$ws.Range("B17").Value = 238.8
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 200.7
$ws.Range("E17").Value = 206.54
$ws.Range("F17").Value = 4.77
$ws.Range("G17").Value = 43.3
$ws.Range("H17").Value = 6.2
$ws.Range("I17").Value = 5.3
$ws.Range("J17").Value = 4.5
$ws.Range("K17").Value = 4.7
$ws.Range("L17").Value = 2.9
$ws.Range("M17").Value = 0.4

# Row 18
$ws.Range("A18").Value = 1974
$ws.Range("B18").Value = 247.3
$ws.Range("C18").Value = 13.6
$ws.Range("D18").Value = 207
$ws.Range("E18").Value = 213.03
$ws.Range("F18").Value = 5.06
$ws.Range("G18").Value = 42.1
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 4.6
$ws.Range("J18").Value = 3.9
$ws.Range("K18").Value = 4.7
$ws.Range("L18").Value = 2.6
$ws.Range("M18").Value = 0.7

# Remove the leftover tail rows now that every column has been compacted
$ws.Range("A19:M34").ClearContents()

# Restore the saved selection
$ws.Range("I6").Select()